# Scheduled runner update: refresh market price snapshot data (columns H:N)
# across the per-job Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 359.375
$ws.Range("I2").Value = 186.36363
$ws.Range("J2").Value = 740
$ws.Range("K2").Value = 186.36363
$ws.Range("L2").Value = 740
$ws.Range("M2").Value = -73.36363
$ws.Range("N2").Value = -966

$ws.Range("H4").Value = 298
$ws.Range("I4").Value = 298
$ws.Range("K4").Value = 298
$ws.Range("M4").Value = -184

$ws.Range("H38").Value = 1976.5714
$ws.Range("I38").Value = 381.5
$ws.Range("J38").Value = 4103.3335
$ws.Range("K38").Value = 1144.5
$ws.Range("L38").Value = 12310.0005
$ws.Range("M38").Value = -772.5
$ws.Range("N38").Value = -13054.0005

$ws.Range("H39").Value = 427
$ws.Range("I39").Value = 104
$ws.Range("J39").Value = 750
$ws.Range("K39").Value = 312
$ws.Range("L39").Value = 2250
$ws.Range("M39").Value = -16
$ws.Range("N39").Value = -2842

$ws.Range("H55").Value = 1498.75
$ws.Range("I55").Value = 1250
$ws.Range("J55").Value = 1747.5
$ws.Range("K55").Value = 1250
$ws.Range("L55").Value = 1747.5
$ws.Range("M55").Value = -1036
$ws.Range("N55").Value = -2175.5

$ws.Range("H58").Value = 2114
$ws.Range("I58").Value = 805
$ws.Range("J58").Value = 3423
$ws.Range("K58").Value = 2415
$ws.Range("L58").Value = 10269
$ws.Range("M58").Value = -2265
$ws.Range("N58").Value = -10569

$ws.Range("H111").Value = 711.75
$ws.Range("I111").Value = 711.75
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 2135.25
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 931.75
$ws.Range("N111").ClearContents()

$ws.Range("H138").Value = 4784.8945
$ws.Range("J138").Value = 4994.9287
$ws.Range("L138").Value = 14984.7861
$ws.Range("N138").Value = -25264.7861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 40.5
$ws.Range("I5").Value = 40.666668
$ws.Range("K5").Value = 40.666668
$ws.Range("M5").Value = 71.333332

$ws.Range("H32").Value = 4856.641
$ws.Range("I32").Value = 4787.1055
$ws.Range("K32").Value = 4787.1055
$ws.Range("M32").Value = -4500.1055

$ws.Range("H46").Value = 14570.5
$ws.Range("J46").Value = 14570.5
$ws.Range("L46").Value = 14570.5
$ws.Range("N46").Value = -15208.5

$ws.Range("H61").Value = 4133.6
$ws.Range("I61").Value = 4556
$ws.Range("K61").Value = 4556
$ws.Range("M61").Value = -4344

$ws.Range("H74").Value = 4713.972
$ws.Range("I74").Value = 4500.3667
$ws.Range("J74").Value = 5782
$ws.Range("K74").Value = 4500.3667
$ws.Range("L74").Value = 5782
$ws.Range("M74").Value = -3626.3667
$ws.Range("N74").Value = -7530

$ws.Range("H77").Value = 4713.972
$ws.Range("I77").Value = 4500.3667
$ws.Range("J77").Value = 5782
$ws.Range("K77").Value = 22501.8335
$ws.Range("L77").Value = 28910
$ws.Range("M77").Value = -18133.8335
$ws.Range("N77").Value = -37646

$ws.Range("H132").Value = 2927.3333
$ws.Range("I132").Value = 2174
$ws.Range("K132").Value = 6522
$ws.Range("M132").Value = -3992

$ws.Range("H136").Value = 4133.6
$ws.Range("I136").Value = 4556
$ws.Range("K136").Value = 13668
$ws.Range("M136").Value = -11118

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 40.5
$ws.Range("I4").Value = 40.666668
$ws.Range("K4").Value = 40.666668
$ws.Range("M4").Value = 74.333332

$ws.Range("H94").Value = 3021.8
$ws.Range("I94").Value = 2777.25
$ws.Range("K94").Value = 2777.25
$ws.Range("M94").Value = -2326.25

$ws.Range("H134").Value = 3507.4736
$ws.Range("I134").Value = 2302.4707
$ws.Range("K134").Value = 6907.4121
$ws.Range("M134").Value = -4372.4121

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1245.1428
$ws.Range("I31").Value = 948.36365
$ws.Range("J31").Value = 2333.3333
$ws.Range("K31").Value = 948.36365
$ws.Range("L31").Value = 2333.3333
$ws.Range("M31").Value = -653.36365
$ws.Range("N31").Value = -2923.3333

$ws.Range("H34").Value = 1245.1428
$ws.Range("I34").Value = 948.36365
$ws.Range("J34").Value = 2333.3333
$ws.Range("K34").Value = 948.36365
$ws.Range("L34").Value = 2333.3333
$ws.Range("M34").Value = -746.36365
$ws.Range("N34").Value = -2737.3333

$ws.Range("H58").Value = 3410.48
$ws.Range("I58").Value = 3357.5908
$ws.Range("J58").Value = 3798.3333
$ws.Range("K58").Value = 3357.5908
$ws.Range("L58").Value = 3798.3333
$ws.Range("M58").Value = -3154.5908
$ws.Range("N58").Value = -4204.3333

$ws.Range("H132").Value = 3651.375
$ws.Range("I132").Value = 2462.6
$ws.Range("K132").Value = 7387.799999999999
$ws.Range("M132").Value = -4857.799999999999

$ws.Range("H134").Value = 6729.8335
$ws.Range("I134").Value = 7596
$ws.Range("K134").Value = 22788
$ws.Range("M134").Value = -20253

$ws.Range("H136").Value = 3410.48
$ws.Range("I136").Value = 3357.5908
$ws.Range("J136").Value = 3798.3333
$ws.Range("K136").Value = 10072.7724
$ws.Range("L136").Value = 11394.9999
$ws.Range("M136").Value = -7522.7724
$ws.Range("N136").Value = -16494.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4699.75
$ws.Range("I81").Value = 3149.5
$ws.Range("K81").Value = 9448.5
$ws.Range("M81").Value = -8325.5

$ws.Range("H84").Value = 4699.75
$ws.Range("I84").Value = 3149.5
$ws.Range("K84").Value = 28345.5
$ws.Range("M84").Value = -22729.5

$ws.Range("H98").Value = 1415
$ws.Range("I98").Value = 800
$ws.Range("J98").Value = 1568.75
$ws.Range("K98").Value = 2400
$ws.Range("L98").Value = 4706.25
$ws.Range("M98").Value = -902
$ws.Range("N98").Value = -7702.25

$ws.Range("H107").Value = 536.6667
$ws.Range("I107").Value = 70
$ws.Range("K107").Value = 210
$ws.Range("M107").Value = 1710

$ws.Range("H121").Value = 129.66667
$ws.Range("I121").Value = 129.66667
$ws.Range("K121").Value = 389.00001
$ws.Range("M121").Value = 920.99999

$ws.Range("H131").Value = 2562.762
$ws.Range("J131").Value = 2481.0588
$ws.Range("L131").Value = 7443.176399999999
$ws.Range("N131").Value = -17523.1764

$ws.Range("H137").Value = 984
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 984
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 2952
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -13152

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 125.8
$ws.Range("I2").Value = 151.14285
$ws.Range("K2").Value = 151.14285
$ws.Range("M2").Value = -38.14285000000001

$ws.Range("H102").Value = 2214.5
$ws.Range("I102").Value = 2214.5
$ws.Range("K102").Value = 2214.5
$ws.Range("M102").Value = -592.5

$ws.Range("H107").Value = 2068.4285
$ws.Range("I107").Value = 1580
$ws.Range("K107").Value = 1580
$ws.Range("M107").Value = 340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 38929.168
$ws.Range("J76").Value = 38929.168
$ws.Range("L76").Value = 38929.168
$ws.Range("N76").Value = -39605.168

$ws.Range("H79").Value = 38929.168
$ws.Range("J79").Value = 38929.168
$ws.Range("L79").Value = 38929.168
$ws.Range("N79").Value = -41269.168

$ws.Range("H93").Value = 862.55554
$ws.Range("I93").Value = 843.9167
$ws.Range("K93").Value = 843.9167
$ws.Range("M93").Value = 404.0833

$ws.Range("H136").Value = 750.6667
$ws.Range("J136").Value = 1224.5
$ws.Range("L136").Value = 3673.5
$ws.Range("N136").Value = -8773.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1934.4584
$ws.Range("I136").Value = 2006.0476
$ws.Range("K136").Value = 6018.142800000001
$ws.Range("M136").Value = -3468.142800000001

Write-Host "Updated market price snapshot cells across all sheets."
